$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "Por favor, preste atenção especial às importações, exportações e declarações. É fundamental evitar a duplicação de declarações para não comprometer o sistema. Verifique se as importações já foram declaradas globalmente ou em algum outro lugar, para que não façamos a mesma declaração duas vezes. Isso é crucial para garantir a integridade do código em JavaScript"
$ws.Range("B7").WrapText = $true

$ws.Range("B7").Select()
